$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2316996
$ws.Range("J17").Value = 2647800.8
$ws.Range("L17").Value = 7943402.399999999
$ws.Range("N17").Value = -7943738.399999999
# Row 58
$ws.Range("H58").Value = 1647.3334
$ws.Range("I58").Value = 1540.75
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 4622.25
$ws.Range("L58").Value = 7500
$ws.Range("M58").Value = -4472.25
$ws.Range("N58").Value = -7800
# Row 62
$ws.Range("H62").Value = 77935.07000000001
$ws.Range("I62").Value = 77935.07000000001
$ws.Range("K62").Value = 77935.07000000001
$ws.Range("M62").Value = -77311.07000000001
# Row 65
$ws.Range("H65").Value = 77935.07000000001
$ws.Range("I65").Value = 77935.07000000001
$ws.Range("K65").Value = 389675.35
$ws.Range("M65").Value = -386555.35
# Row 70
$ws.Range("H70").Value = 205910
$ws.Range("I70").Value = 3400
$ws.Range("J70").Value = 228411.11
$ws.Range("K70").Value = 10200
$ws.Range("L70").Value = 685233.33
$ws.Range("M70").Value = -9930
$ws.Range("N70").Value = -685773.33
# Row 73
$ws.Range("H73").Value = 205910
$ws.Range("I73").Value = 3400
$ws.Range("J73").Value = 228411.11
$ws.Range("K73").Value = 10200
$ws.Range("L73").Value = 685233.33
$ws.Range("M73").Value = -9264
$ws.Range("N73").Value = -687105.33
# Row 96
$ws.Range("H96").Value = 189.125
$ws.Range("J96").Value = 215.5
$ws.Range("L96").Value = 646.5
$ws.Range("N96").Value = -3392.5
# Row 116
$ws.Range("H116").Value = 5998.4443
$ws.Range("I116").Value = 6664.8335
$ws.Range("K116").Value = 6664.8335
$ws.Range("M116").Value = -3222.8335
# Row 132
$ws.Range("H132").Value = 1212.2307
$ws.Range("I132").Value = 1129.1082
$ws.Range("K132").Value = 3387.3246
$ws.Range("M132").Value = -857.3245999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 5234.7354
$ws.Range("I61").Value = 4331.222
$ws.Range("K61").Value = 4331.222
$ws.Range("M61").Value = -4119.222
# Row 136
$ws.Range("H136").Value = 5234.7354
$ws.Range("I136").Value = 4331.222
$ws.Range("K136").Value = 12993.666
$ws.Range("M136").Value = -10443.666

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1245.1
$ws.Range("I107").Value = 1159.7778
$ws.Range("K107").Value = 1159.7778
$ws.Range("M107").Value = 760.2221999999999
# Row 134
$ws.Range("H134").Value = 2254.5874
$ws.Range("I134").Value = 2175.2712
$ws.Range("J134").Value = 3424.5
$ws.Range("K134").Value = 6525.8136
$ws.Range("L134").Value = 10273.5
$ws.Range("M134").Value = -3990.8136
$ws.Range("N134").Value = -15343.5
# Row 135
$ws.Range("H135").Value = 74437.14
$ws.Range("J135").Value = 74437.14
$ws.Range("L135").Value = 74437.14
$ws.Range("N135").Value = -84577.14

$ws = $wb.Worksheets.Item("CRP")
# Row 93
$ws.Range("H93").Value = 10637.857
$ws.Range("I93").Value = 5910.8335
$ws.Range("K93").Value = 5910.8335
$ws.Range("M93").Value = -4038.8335
# Row 94
$ws.Range("H94").Value = 3119
$ws.Range("J94").Value = 3360
$ws.Range("L94").Value = 3360
$ws.Range("N94").Value = -4262
# Row 99
$ws.Range("H99").Value = 2874.4
$ws.Range("I99").Value = 2756
$ws.Range("J99").Value = 3348
$ws.Range("K99").Value = 2756
$ws.Range("L99").Value = 3348
$ws.Range("M99").Value = -1258
$ws.Range("N99").Value = -6344
# Row 103
$ws.Range("H103").Value = 28333.166
$ws.Range("I103").Value = 25999.8
$ws.Range("K103").Value = 25999.8
$ws.Range("M103").Value = -24827.8
# Row 122
$ws.Range("H122").Value = 1964.2222
$ws.Range("I122").Value = 1545.909
$ws.Range("K122").Value = 4637.727000000001
$ws.Range("M122").Value = -2187.727000000001
# Row 126
$ws.Range("H126").Value = 2874.4
$ws.Range("I126").Value = 2756
$ws.Range("J126").Value = 3348
$ws.Range("K126").Value = 8268
$ws.Range("L126").Value = 10044
$ws.Range("M126").Value = -5798
$ws.Range("N126").Value = -14984
# Row 134
$ws.Range("H134").Value = 10292.415
$ws.Range("I134").Value = 6150.5557
$ws.Range("K134").Value = 18451.6671
$ws.Range("M134").Value = -15916.6671

$ws = $wb.Worksheets.Item("CUL")
# Row 106
$ws.Range("H106").Value = 7750
$ws.Range("J106").Value = 7333.3335
$ws.Range("L106").Value = 22000.0005
$ws.Range("N106").Value = -23892.0005
# Row 121
$ws.Range("H121").Value = 25532.5
$ws.Range("J121").Value = 1000
$ws.Range("L121").Value = 3000
$ws.Range("N121").Value = -5620
# Row 122
$ws.Range("H122").Value = 1074.238
$ws.Range("I122").Value = 457
$ws.Range("J122").Value = 1382.8572
$ws.Range("K122").Value = 4113
$ws.Range("L122").Value = 12445.7148
$ws.Range("M122").Value = -1663
$ws.Range("N122").Value = -17345.7148
# Row 140
$ws.Range("H140").Value = 3974.7
$ws.Range("I140").Value = 3305.2222
$ws.Range("K140").Value = 9915.6666
$ws.Range("M140").Value = -4735.6666

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 21183.166
$ws.Range("I126").Value = 24549.9
$ws.Range("K126").Value = 73649.70000000001
$ws.Range("M126").Value = -71179.70000000001
# Row 132
$ws.Range("H132").Value = 5625.125
$ws.Range("I132").Value = 4748.75
$ws.Range("J132").Value = 6501.5
$ws.Range("K132").Value = 14246.25
$ws.Range("L132").Value = 19504.5
$ws.Range("M132").Value = -11716.25
$ws.Range("N132").Value = -24564.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7948.421
$ws.Range("I7").Value = 7813.25
$ws.Range("K7").Value = 7813.25
$ws.Range("M7").Value = -7701.25
# Row 126
$ws.Range("H126").Value = 7948.421
$ws.Range("I126").Value = 7813.25
$ws.Range("K126").Value = 23439.75
$ws.Range("M126").Value = -20969.75

$ws = $wb.Worksheets.Item("WVR")
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").ClearContents()
$ws.Range("N59").ClearContents()
# Row 70
$ws.Range("H70").Value = 52550
$ws.Range("J70").Value = 52550
$ws.Range("L70").Value = 52550
$ws.Range("N70").Value = -53180
# Row 73
$ws.Range("H73").Value = 52550
$ws.Range("J73").Value = 52550
$ws.Range("L73").Value = 52550
$ws.Range("N73").Value = -54734
# Row 132
$ws.Range("H132").Value = 2211.6135
$ws.Range("I132").Value = 2008.1842
$ws.Range("K132").Value = 6024.5526
$ws.Range("M132").Value = -3494.5526
# Row 136
$ws.Range("H136").Value = 1942.65
$ws.Range("I136").Value = 1657.1333
$ws.Range("J136").Value = 2799.2
$ws.Range("K136").Value = 4971.3999
$ws.Range("L136").Value = 8397.599999999999
$ws.Range("M136").Value = -2421.3999
$ws.Range("N136").Value = -13497.6
